# Scheduled-runner market price refresh for Sheets workbook
# (values sourced from an external price lookup; no formulas involved)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 787.7646999999999  # H2: 608.087 -> 787.7646999999999
$ws.Cells.Item(2, 10).Value = 1237  # J2: 781.8 -> 1237
$ws.Cells.Item(2, 12).Value = 1237  # L2: 781.8 -> 1237
$ws.Cells.Item(2, 14).Value = -1463  # N2: -1007.8 -> -1463
$ws.Cells.Item(98, 8).Value = 975.2  # H98: 1077.4445 -> 975.2
$ws.Cells.Item(98, 9).Value = 1158.8334  # I98: 1379.6 -> 1158.8334
$ws.Cells.Item(98, 11).Value = 1158.8334  # K98: 1379.6 -> 1158.8334
$ws.Cells.Item(98, 13).Value = 339.1666  # M98: 118.4000000000001 -> 339.1666
$ws.Cells.Item(122, 8).Value = 975.2  # H122: 1077.4445 -> 975.2
$ws.Cells.Item(122, 9).Value = 1158.8334  # I122: 1379.6 -> 1158.8334
$ws.Cells.Item(122, 11).Value = 3476.5002  # K122: 4138.799999999999 -> 3476.5002
$ws.Cells.Item(122, 13).Value = -1026.5002  # M122: -1688.799999999999 -> -1026.5002
$ws.Cells.Item(132, 8).Value = 3369.125  # H132: 3452.8572 -> 3369.125
$ws.Cells.Item(132, 9).Value = 3421.8572  # I132: 3528.3333 -> 3421.8572
$ws.Cells.Item(132, 11).Value = 10265.5716  # K132: 10584.9999 -> 10265.5716
$ws.Cells.Item(132, 13).Value = -7735.571599999999  # M132: -8054.999899999999 -> -7735.571599999999
$ws.Cells.Item(138, 8).Value = 4243.25  # H138: 4801.25 -> 4243.25
$ws.Cells.Item(138, 9).Value = 1998.5  # I138: 2262.5715 -> 1998.5
$ws.Cells.Item(138, 10).Value = 4804.4375  # J138: 5647.476 -> 4804.4375
$ws.Cells.Item(138, 11).Value = 5995.5  # K138: 6787.7145 -> 5995.5
$ws.Cells.Item(138, 12).Value = 14413.3125  # L138: 16942.428 -> 14413.3125
$ws.Cells.Item(138, 13).Value = -855.5  # M138: -1647.7145 -> -855.5
$ws.Cells.Item(138, 14).Value = -24693.3125  # N138: -27222.428 -> -24693.3125

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 1167.375  # H74: 2089.7144 -> 1167.375
$ws.Cells.Item(74, 9).Value = 991.4286  # I74: 2089.7144 -> 991.4286
$ws.Cells.Item(74, 10).Value = 2399  # J74: 0 -> 2399
$ws.Cells.Item(74, 11).Value = 991.4286  # K74: 2089.7144 -> 991.4286
$ws.Cells.Item(74, 12).Value = 2399  # L74: 0 -> 2399
$ws.Cells.Item(74, 13).Value = -117.4286  # M74: -1215.7144 -> -117.4286
$ws.Cells.Item(74, 14).Value = -4147  # N74: (none) -> -4147
$ws.Cells.Item(77, 8).Value = 1167.375  # H77: 2089.7144 -> 1167.375
$ws.Cells.Item(77, 9).Value = 991.4286  # I77: 2089.7144 -> 991.4286
$ws.Cells.Item(77, 10).Value = 2399  # J77: 0 -> 2399
$ws.Cells.Item(77, 11).Value = 4957.143  # K77: 10448.572 -> 4957.143
$ws.Cells.Item(77, 12).Value = 11995  # L77: 0 -> 11995
$ws.Cells.Item(77, 13).Value = -589.143  # M77: -6080.572 -> -589.143
$ws.Cells.Item(77, 14).Value = -20731  # N77: (none) -> -20731

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2181.9167  # H20: 2300.2727 -> 2181.9167
$ws.Cells.Item(20, 10).Value = 2712.1667  # J20: 3078.6 -> 2712.1667
$ws.Cells.Item(20, 12).Value = 2712.1667  # L20: 3078.6 -> 2712.1667
$ws.Cells.Item(20, 14).Value = -3206.1667  # N20: -3572.6 -> -3206.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2035.9  # H31: 1742 -> 2035.9
$ws.Cells.Item(31, 9).Value = 1769.8889  # I31: 1742 -> 1769.8889
$ws.Cells.Item(31, 10).Value = 4430  # J31: 0 -> 4430
$ws.Cells.Item(31, 11).Value = 1769.8889  # K31: 1742 -> 1769.8889
$ws.Cells.Item(31, 12).Value = 4430  # L31: 0 -> 4430
$ws.Cells.Item(31, 13).Value = -1474.8889  # M31: -1447 -> -1474.8889
$ws.Cells.Item(31, 14).Value = -5020  # N31: (none) -> -5020
$ws.Cells.Item(34, 8).Value = 2035.9  # H34: 1742 -> 2035.9
$ws.Cells.Item(34, 9).Value = 1769.8889  # I34: 1742 -> 1769.8889
$ws.Cells.Item(34, 10).Value = 4430  # J34: 0 -> 4430
$ws.Cells.Item(34, 11).Value = 1769.8889  # K34: 1742 -> 1769.8889
$ws.Cells.Item(34, 12).Value = 4430  # L34: 0 -> 4430
$ws.Cells.Item(34, 13).Value = -1567.8889  # M34: -1540 -> -1567.8889
$ws.Cells.Item(34, 14).Value = -4834  # N34: (none) -> -4834
$ws.Cells.Item(58, 8).Value = 2393.5  # H58: 2462.1428 -> 2393.5
$ws.Cells.Item(58, 9).Value = 1980.2  # I58: 1997 -> 1980.2
$ws.Cells.Item(58, 11).Value = 1980.2  # K58: 1997 -> 1980.2
$ws.Cells.Item(58, 13).Value = -1777.2  # M58: -1794 -> -1777.2
$ws.Cells.Item(70, 8).Value = 35000  # H70: 0 -> 35000
$ws.Cells.Item(70, 10).Value = 35000  # J70: 0 -> 35000
$ws.Cells.Item(70, 12).Value = 35000  # L70: 0 -> 35000
$ws.Cells.Item(70, 14).Value = -35630  # N70: (none) -> -35630
$ws.Cells.Item(73, 8).Value = 35000  # H73: 0 -> 35000
$ws.Cells.Item(73, 10).Value = 35000  # J73: 0 -> 35000
$ws.Cells.Item(73, 12).Value = 35000  # L73: 0 -> 35000
$ws.Cells.Item(73, 14).Value = -37184  # N73: (none) -> -37184
$ws.Cells.Item(86, 8).Value = 11055.625  # H86: 12049.444 -> 11055.625
$ws.Cells.Item(86, 10).Value = 6993  # J86: 11328.667 -> 6993
$ws.Cells.Item(86, 12).Value = 6993  # L86: 11328.667 -> 6993
$ws.Cells.Item(86, 14).Value = -9239  # N86: -13574.667 -> -9239
$ws.Cells.Item(89, 8).Value = 11055.625  # H89: 12049.444 -> 11055.625
$ws.Cells.Item(89, 10).Value = 6993  # J89: 11328.667 -> 6993
$ws.Cells.Item(89, 12).Value = 34965  # L89: 56643.335 -> 34965
$ws.Cells.Item(89, 14).Value = -46197  # N89: -67875.33499999999 -> -46197
$ws.Cells.Item(94, 8).Value = 5474  # H94: 3443.2144 -> 5474
$ws.Cells.Item(94, 9).Value = 11999  # I94: 2551.5557 -> 11999
$ws.Cells.Item(94, 10).Value = 4169  # J94: 5048.2 -> 4169
$ws.Cells.Item(94, 11).Value = 11999  # K94: 2551.5557 -> 11999
$ws.Cells.Item(94, 12).Value = 4169  # L94: 5048.2 -> 4169
$ws.Cells.Item(94, 13).Value = -11548  # M94: -2100.5557 -> -11548
$ws.Cells.Item(94, 14).Value = -5071  # N94: -5950.2 -> -5071
$ws.Cells.Item(132, 8).Value = 3979.6667  # H132: 4519.75 -> 3979.6667
$ws.Cells.Item(132, 9).Value = 4960.6665  # I132: 5991.5 -> 4960.6665
$ws.Cells.Item(132, 10).Value = 2998.6667  # J132: 3048 -> 2998.6667
$ws.Cells.Item(132, 11).Value = 14881.9995  # K132: 17974.5 -> 14881.9995
$ws.Cells.Item(132, 12).Value = 8996.000100000001  # L132: 9144 -> 8996.000100000001
$ws.Cells.Item(132, 13).Value = -12351.9995  # M132: -15444.5 -> -12351.9995
$ws.Cells.Item(132, 14).Value = -14056.0001  # N132: -14204 -> -14056.0001
$ws.Cells.Item(134, 8).Value = 2340.1333  # H134: 2365.9285 -> 2340.1333
$ws.Cells.Item(134, 9).Value = 2348.5715  # I134: 2377 -> 2348.5715
$ws.Cells.Item(134, 11).Value = 7045.7145  # K134: 7131 -> 7045.7145
$ws.Cells.Item(134, 13).Value = -4510.7145  # M134: -4596 -> -4510.7145
$ws.Cells.Item(136, 8).Value = 2393.5  # H136: 2462.1428 -> 2393.5
$ws.Cells.Item(136, 9).Value = 1980.2  # I136: 1997 -> 1980.2
$ws.Cells.Item(136, 11).Value = 5940.6  # K136: 5991 -> 5940.6
$ws.Cells.Item(136, 13).Value = -3390.6  # M136: -3441 -> -3390.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 2360  # H34: 483 -> 2360
$ws.Cells.Item(34, 10).Value = 4245.5  # J34: 500 -> 4245.5
$ws.Cells.Item(34, 12).Value = 12736.5  # L34: 1500 -> 12736.5
$ws.Cells.Item(34, 14).Value = -12904.5  # N34: -1668 -> -12904.5
$ws.Cells.Item(55, 8).Value = 3245  # H55: 500 -> 3245
$ws.Cells.Item(55, 10).Value = 5990  # J55: 0 -> 5990
$ws.Cells.Item(55, 12).Value = 17970  # L55: 0 -> 17970
$ws.Cells.Item(55, 14).Value = -18324  # N55: (none) -> -18324
$ws.Cells.Item(92, 8).Value = 849.5  # H92: 473.5 -> 849.5
$ws.Cells.Item(92, 9).Value = 700  # I92: 400 -> 700
$ws.Cells.Item(92, 10).Value = 999  # J92: 547 -> 999
$ws.Cells.Item(92, 11).Value = 2100  # K92: 1200 -> 2100
$ws.Cells.Item(92, 12).Value = 2997  # L92: 1641 -> 2997
$ws.Cells.Item(92, 13).Value = -852  # M92: 48 -> -852
$ws.Cells.Item(92, 14).Value = -5493  # N92: -4137 -> -5493
$ws.Cells.Item(114, 8).Value = 3010  # H114: 3301.2856 -> 3010
$ws.Cells.Item(114, 9).Value = 1000  # I114: 2014 -> 1000
$ws.Cells.Item(114, 10).Value = 4015  # J114: 3816.2 -> 4015
$ws.Cells.Item(114, 11).Value = 3000  # K114: 6042 -> 3000
$ws.Cells.Item(114, 12).Value = 12045  # L114: 11448.6 -> 12045
$ws.Cells.Item(114, 13).Value = 254  # M114: -2788 -> 254
$ws.Cells.Item(114, 14).Value = -18553  # N114: -17956.6 -> -18553

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(35, 8).Value = 0  # H35: 100000 -> 0
$ws.Cells.Item(35, 9).Value = 0  # I35: 100000 -> 0
$ws.Cells.Item(35, 11).Value = 0  # K35: 100000 -> 0
$ws.Cells.Item(35, 13).ClearContents()  # M35: -99702 -> (removed)
$ws.Cells.Item(132, 8).Value = 13047.3  # H132: 14821.75 -> 13047.3
$ws.Cells.Item(132, 9).Value = 13730.444  # I132: 14821.75 -> 13730.444
$ws.Cells.Item(132, 10).Value = 6899  # J132: 0 -> 6899
$ws.Cells.Item(132, 11).Value = 41191.33199999999  # K132: 44465.25 -> 41191.33199999999
$ws.Cells.Item(132, 12).Value = 20697  # L132: 0 -> 20697
$ws.Cells.Item(132, 13).Value = -38661.33199999999  # M132: -41935.25 -> -38661.33199999999
$ws.Cells.Item(132, 14).Value = -25757  # N132: (none) -> -25757

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 10601.333  # H132: 9306.75 -> 10601.333
$ws.Cells.Item(132, 9).Value = 11961.6  # I132: 9779.143 -> 11961.6
$ws.Cells.Item(132, 10).Value = 3800  # J132: 6000 -> 3800
$ws.Cells.Item(132, 11).Value = 35884.8  # K132: 29337.429 -> 35884.8
$ws.Cells.Item(132, 12).Value = 11400  # L132: 18000 -> 11400
$ws.Cells.Item(132, 13).Value = -33354.8  # M132: -26807.429 -> -33354.8
$ws.Cells.Item(132, 14).Value = -16460  # N132: -23060 -> -16460
$ws.Cells.Item(137, 8).Value = 102800.336  # H137: 111631.836 -> 102800.336
$ws.Cells.Item(137, 9).Value = 98645  # I137: 99340 -> 98645
$ws.Cells.Item(137, 10).Value = 111111  # J137: 117777.75 -> 111111
$ws.Cells.Item(137, 11).Value = 98645  # K137: 99340 -> 98645
$ws.Cells.Item(137, 12).Value = 111111  # L137: 117777.75 -> 111111
$ws.Cells.Item(137, 13).Value = -93545  # M137: -94240 -> -93545
$ws.Cells.Item(137, 14).Value = -121311  # N137: -127977.75 -> -121311

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 19648.857  # H41: 19698.166 -> 19648.857
$ws.Cells.Item(41, 9).Value = 19567.666  # I41: 19672.5 -> 19567.666
$ws.Cells.Item(41, 10).Value = 19709.75  # J41: 19711 -> 19709.75
$ws.Cells.Item(41, 11).Value = 19567.666  # K41: 19672.5 -> 19567.666
$ws.Cells.Item(41, 12).Value = 19709.75  # L41: 19711 -> 19709.75
$ws.Cells.Item(41, 13).Value = -19177.666  # M41: -19282.5 -> -19177.666
$ws.Cells.Item(41, 14).Value = -20489.75  # N41: -20491 -> -20489.75
$ws.Cells.Item(45, 8).Value = 21150  # H45: 21133 -> 21150
$ws.Cells.Item(45, 10).Value = 21150  # J45: 21133 -> 21150
$ws.Cells.Item(45, 12).Value = 21150  # L45: 21133 -> 21150
$ws.Cells.Item(45, 14).Value = -22132  # N45: -22115 -> -22132
$ws.Cells.Item(101, 8).Value = 15000  # H101: 16151 -> 15000
$ws.Cells.Item(101, 10).Value = 15000  # J101: 16151 -> 15000
$ws.Cells.Item(101, 12).Value = 15000  # L101: 16151 -> 15000
$ws.Cells.Item(101, 14).Value = -21490  # N101: -22641 -> -21490
$ws.Cells.Item(132, 8).Value = 2606.4285  # H132: 2281.25 -> 2606.4285
$ws.Cells.Item(132, 10).Value = 797  # J132: 401 -> 797
$ws.Cells.Item(132, 12).Value = 2391  # L132: 1203 -> 2391
$ws.Cells.Item(132, 14).Value = -7451  # N132: -6263 -> -7451
$ws.Cells.Item(136, 8).Value = 2467.2666  # H136: 2363 -> 2467.2666
$ws.Cells.Item(136, 9).Value = 2777.0833  # I136: 2624.923 -> 2777.0833
$ws.Cells.Item(136, 11).Value = 8331.249899999999  # K136: 7874.768999999999 -> 8331.249899999999
$ws.Cells.Item(136, 13).Value = -5781.249899999999  # M136: -5324.768999999999 -> -5781.249899999999
$ws.Cells.Item(137, 8).Value = 0  # H137: 43000 -> 0
$ws.Cells.Item(137, 10).Value = 0  # J137: 43000 -> 0
$ws.Cells.Item(137, 12).Value = 0  # L137: 43000 -> 0
$ws.Cells.Item(137, 14).ClearContents()  # N137: -53200 -> (removed)
$ws.Cells.Item(138, 8).Value = 79999  # H138: 0 -> 79999
$ws.Cells.Item(138, 10).Value = 79999  # J138: 0 -> 79999
$ws.Cells.Item(138, 12).Value = 79999  # L138: 0 -> 79999
$ws.Cells.Item(138, 14).Value = -90279  # N138: (none) -> -90279
$ws.Cells.Item(140, 8).Value = 89995  # H140: 89994.5 -> 89995
$ws.Cells.Item(140, 10).Value = 89995  # J140: 89994.5 -> 89995
$ws.Cells.Item(140, 12).Value = 89995  # L140: 89994.5 -> 89995
$ws.Cells.Item(140, 14).Value = -100355  # N140: -100354.5 -> -100355
$ws.Cells.Item(141, 8).Value = 297999.5  # H141: 231999.67 -> 297999.5
$ws.Cells.Item(141, 10).Value = 297999.5  # J141: 231999.67 -> 297999.5
$ws.Cells.Item(141, 12).Value = 297999.5  # L141: 231999.67 -> 297999.5
$ws.Cells.Item(141, 14).Value = -308359.5  # N141: -242359.67 -> -308359.5
